$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")
$ws.Activate()
$ws.Range("B3").Value = "Nancy@1234"
$ws.Range("B3").Select()
